$wb = $excel.ActiveWorkbook
$wsTemplate = $wb.Worksheets.Item("Template")
$wsOptions  = $wb.Worksheets.Item("Options")

# =========================================================================
# Template sheet: shift the tail of the header row right by one column and
# insert the new "Sample_type" header in the vacated column P.
#   old P1 Sample_or_Control -> Q1
#   old Q1 SequencingRun     -> R1
#   old R1 Notes             -> S1
#   new P1                   -> Sample_type
# =========================================================================
$wsTemplate.Range("S1").Value = $wsTemplate.Range("R1").Value()
$wsTemplate.Range("R1").Value = $wsTemplate.Range("Q1").Value()
$wsTemplate.Range("Q1").Value = $wsTemplate.Range("P1").Value()
$wsTemplate.Range("P1").Value = "Sample_type"

# =========================================================================
# Options sheet: shift the "Sample_or_Control" option list (old P1:P2) right
# into Q1:Q2, then fill column P with the new Sample_type options.
# =========================================================================
$wsOptions.Range("Q1").Value = $wsOptions.Range("P1").Value()
$wsOptions.Range("Q2").Value = $wsOptions.Range("P2").Value()
$wsOptions.Range("P1").Value = "soil"
$wsOptions.Range("P2").Value = "root"
$wsOptions.Range("P3").Value = "marine_sediment"
$wsOptions.Range("P4").Value = "marine_root"

# New "NA" choices appended to several existing option lists.
$wsOptions.Range("H3").Value  = "NA"
$wsOptions.Range("L4").Value  = "NA"
$wsOptions.Range("M5").Value  = "NA"
$wsOptions.Range("G12").Value = "NA"
$wsOptions.Range("I15").Value = "NA"

# =========================================================================
# Template sheet data validations: grow the existing lists by one row each,
# repoint the existing "P" rule (now Sample_type) at the new soil/root/... list,
# and add a brand-new rule for the shifted "Sample_or_Control" column (Q).
# =========================================================================
$wsTemplate.Range("G2:G101").Validation.Formula1 = "=Options!`$G`$1:`$G`$12"
$wsTemplate.Range("H2:H101").Validation.Formula1 = "=Options!`$H`$1:`$H`$3"
$wsTemplate.Range("I2:I101").Validation.Formula1 = "=Options!`$I`$1:`$I`$15"
$wsTemplate.Range("L2:L101").Validation.Formula1 = "=Options!`$L`$1:`$L`$4"
$wsTemplate.Range("M2:M101").Validation.Formula1 = "=Options!`$M`$1:`$M`$5"
$wsTemplate.Range("P2:P101").Validation.Formula1 = "=Options!`$P`$1:`$P`$4"

$wsTemplate.Range("Q2:Q101").Validation.Add(3, 1, 1, "=Options!`$Q`$1:`$Q`$2")
